$d = $word.ActiveDocument

# 1) "Реализовать веб-сервис из задачи 7" -> "Реализовать алгоритм из задачи 7"
$d.Content.Find.Execute(
    "веб-сервис из задачи 7",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "алгоритм из задачи 7",
    2
)

# 2) Trim the trailing clause and replace the ending of task 8's paragraph
$d.Content.Find.Execute(
    "запроса, а не над несколькими независимыми запросами. Предусмотреть корректное завершение сервиса без использования консоли (она плохо работает в MPI.NET).",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "запроса. Предусмотреть корректное завершение работы отдельных процессов.",
    2
)
